$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the set of cell -> new text value updates described by the diff.
$updates = @{
    "D2" = "330.90"
    "E2" = "0.02%"
    "D3" = "45.34"
    "E3" = "2.35%"
    "D4" = "5.586"
    "E4" = "2.04%"
    "D5" = "0.08336"
    "E5" = "3.78%"
    "D6" = "2.113"
    "E6" = "6.25%"
    "D7" = "0.9759"
    "E7" = "2.36%"
    "E8" = "-0.60%"
    "D9" = "0.1196"
    "E9" = "4.95%"
    "D10" = "0.1923"
    "E10" = "1.06%"
    "D11" = "10.32"
    "E11" = "-2.99%"
    "D12" = "0.09820"
    "E12" = "-0.86%"
    "D13" = "0.04679"
    "E13" = "-3.21%"
    "D14" = "0.1058"
    "E14" = "-0.43%"
    "D15" = "0.001297"
    "E15" = "3.18%"
    "D16" = "0.006053"
    "E16" = "1.07%"
    "D18" = "4.449"
    "E18" = "1.34%"
    "D19" = "0.3342"
    "E19" = "-2.50%"
    "D20" = "0.1392"
    "E20" = "-0.41%"
    "D21" = "0.2674"
    "E21" = "6.96%"
    "D22" = "0.04163"
    "E22" = "2.24%"
    "D23" = "0.001293"
    "E23" = "1.70%"
    "D24" = "0.004580"
    "E24" = "4.88%"
    "D25" = "0.0001303"
    "E25" = "8.76%"
    "D26" = "0.0003748"
    "E26" = "0.28%"
    "D38" = "0.02704"
    "E38" = "3.76%"
    "D39" = "0.05742"
    "E39" = "-1.39%"
    "D40" = "0.007855"
    "E40" = "3.87%"
    "D41" = "0.1433"
    "E41" = "1.89%"
    "D42" = "0.007471"
    "E42" = "1.74%"
    "D43" = "0.002100"
    "E43" = "4.42%"
    "D44" = "0.008522"
    "E44" = "-3.74%"
    "D45" = "0.3373"
    "D46" = "0.00007117"
    "E46" = "1.74%"
    "E47" = "0.39%"
    "E48" = "0.47%"
    "D49" = "0.003527"
    "E49" = "0.98%"
    "D50" = "0.003506"
    "E50" = "-0.52%"
    "D51" = "0.00002105"
    "E51" = "0.39%"
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    # Preserve text storage (matches original inlineStr cells) so Excel
    # does not auto-convert numeric-looking strings / percentages into
    # actual numbers.
    $range.NumberFormat = "@"
    $range.Value = $updates[$cell]
}
